$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
}

Set-TextValue $ws.Range("D2") "22.411.03"
$ws.Range("E2").Value = "  -0.12%  "
Set-TextValue $ws.Range("D3") "1.563.02"
$ws.Range("E3").Value = "  -0.63%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -0.26%  "
Set-TextValue $ws.Range("D5") "1.000"
$ws.Range("E5").Value = "  -0.22%  "
Set-TextValue $ws.Range("D6") "284.51"
$ws.Range("E6").Value = "  -2.45%  "
Set-TextValue $ws.Range("D7") "0.3617"
$ws.Range("E7").Value = "  -3.13%  "
Set-TextValue $ws.Range("D8") "48.35"
$ws.Range("E8").Value = "  -2.93%  "
Set-TextValue $ws.Range("D9") "0.3330"
$ws.Range("E9").Value = "  -1.85%  "
Set-TextValue $ws.Range("D10") "1.125"
$ws.Range("E10").Value = "  -0.98%  "
Set-TextValue $ws.Range("D11") "0.07384"
$ws.Range("E11").Value = "  -2.24%  "
Set-TextValue $ws.Range("D12") "1.000"
$ws.Range("E12").Value = "  -0.26%  "
Set-TextValue $ws.Range("D13") "20.71"
$ws.Range("E13").Value = "  -3.14%  "
Set-TextValue $ws.Range("D14") "5.920"
$ws.Range("E14").Value = "  -1.18%  "
Set-TextValue $ws.Range("D15") "6.880"
$ws.Range("E15").Value = "  -0.70%  "
Set-TextValue $ws.Range("D16") "1.562.14"
$ws.Range("E16").Value = "  -0.47%  "
Set-TextValue $ws.Range("D17") "0.00001102"
$ws.Range("E17").Value = "  -1.70%  "
Set-TextValue $ws.Range("D18") "87.91"
$ws.Range("E18").Value = "  -3.39%  "
Set-TextValue $ws.Range("D19") "0.06691"
$ws.Range("E19").Value = "  -0.65%  "
Set-TextValue $ws.Range("D20") "1.000"
$ws.Range("E20").Value = "  -0.21%  "
Set-TextValue $ws.Range("D21") "6.333"
$ws.Range("E21").Value = "  +0.79%  "
Set-TextValue $ws.Range("D22") "16.07"
$ws.Range("E22").Value = "  -1.80%  "
Set-TextValue $ws.Range("D23") "11.98"
$ws.Range("E23").Value = "  -1.37%  "
Set-TextValue $ws.Range("D24") "22.403.80"
$ws.Range("E24").Value = "  -0.13%  "
Set-TextValue $ws.Range("D25") "2.417"
$ws.Range("E25").Value = "  +3.55%  "
Set-TextValue $ws.Range("D26") "2.541"
$ws.Range("E26").Value = "  -3.21%  "
Set-TextValue $ws.Range("D27") "149.94"
$ws.Range("E27").Value = "  +1.06%  "
Set-TextValue $ws.Range("D28") "19.34"
$ws.Range("E28").Value = "  -3.73%  "
Set-TextValue $ws.Range("D29") "5.006"
$ws.Range("E29").Value = "  -0.20%  "
Set-TextValue $ws.Range("D30") "123.02"
$ws.Range("E30").Value = "  -2.08%  "
Set-TextValue $ws.Range("D31") "1.736.59"
$ws.Range("E31").Value = "  -0.63%  "
Set-TextValue $ws.Range("D32") "1.052"
$ws.Range("E32").Value = "  +0.32%  "
Set-TextValue $ws.Range("D33") "6.110"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("E34").Value = "  +0.78%  "
Set-TextValue $ws.Range("D35") "9.768"
$ws.Range("E35").Value = "  -0.10%  "
Set-TextValue $ws.Range("D36") "0.08279"
$ws.Range("E36").Value = "  -0.95%  "
Set-TextValue $ws.Range("D37") "0.02394"
$ws.Range("E37").Value = "  -2.80%  "
Set-TextValue $ws.Range("D38") "0.06381"
$ws.Range("E38").Value = "  -2.05%  "
Set-TextValue $ws.Range("D39") "0.2205"
$ws.Range("E39").Value = "  -3.56%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D40") "1.284"
$ws.Range("E40").Value = "  -7.29%  "
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D41") "5.304"
$ws.Range("E41").Value = "  -2.85%  "
Set-TextValue $ws.Range("D42") "11.11"
$ws.Range("E42").Value = "  -1.63%  "
Set-TextValue $ws.Range("D43") "0.6052"
$ws.Range("E43").Value = "  -2.67%  "
Set-TextValue $ws.Range("D44") "0.9999"
$ws.Range("E44").Value = "  -0.25%  "
Set-TextValue $ws.Range("D45") "13.78"
$ws.Range("E45").Value = "  -0.60%  "
Set-TextValue $ws.Range("D46") "3.749"
$ws.Range("E46").Value = "  -1.66%  "
Set-TextValue $ws.Range("D47") "0.5746"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D48") "2.007"
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D49") "124.34"
$ws.Range("E49").Value = "  -4.07%  "
Set-TextValue $ws.Range("D50") "1.212"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  -1.69%  "
